$wb = $excel.ActiveWorkbook

# Updated "想去人数" (want-to-go count) values for rows 2-7, column F
$updates = @{
    "F2" = 1285
    "F3" = 1640
    "F4" = 64
    "F5" = 6212
    "F6" = 35
    "F7" = 105
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
